$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5155.826
$ws.Range("I40").Value = 4837.5713
$ws.Range("K40").Value = 4837.5713
$ws.Range("M40").Value = -4662.5713
$ws.Range("H138").Value = 1659.62
$ws.Range("I138").Value = 1175.0834
$ws.Range("K138").Value = 3525.2502
$ws.Range("M138").Value = 1614.7498

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2982.139
$ws.Range("I61").Value = 2917.0625
$ws.Range("J61").Value = 3502.75
$ws.Range("K61").Value = 2917.0625
$ws.Range("L61").Value = 3502.75
$ws.Range("M61").Value = -2705.0625
$ws.Range("N61").Value = -3926.75
$ws.Range("H80").Value = 43749.5
$ws.Range("I80").Value = 37500
$ws.Range("K80").Value = 37500
$ws.Range("M80").Value = -36502
$ws.Range("H83").Value = 43749.5
$ws.Range("I83").Value = 37500
$ws.Range("K83").Value = 112500
$ws.Range("M83").Value = -107508
$ws.Range("H122").Value = 4579.0835
$ws.Range("I122").Value = 4210
$ws.Range("J122").Value = 4891.385
$ws.Range("K122").Value = 12630
$ws.Range("L122").Value = 14674.155
$ws.Range("M122").Value = -10180
$ws.Range("N122").Value = -19574.155
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H136").Value = 2982.139
$ws.Range("I136").Value = 2917.0625
$ws.Range("J136").Value = 3502.75
$ws.Range("K136").Value = 8751.1875
$ws.Range("L136").Value = 10508.25
$ws.Range("M136").Value = -6201.1875
$ws.Range("N136").Value = -15608.25
$ws.Range("H137").Value = 499999
$ws.Range("J137").Value = 499999
$ws.Range("L137").Value = 499999
$ws.Range("N137").Value = -510199

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2751.6
$ws.Range("I20").Value = 3384.5386
$ws.Range("K20").Value = 3384.5386
$ws.Range("M20").Value = -3137.5386

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 37.785713
$ws.Range("I7").Value = 49.166668
$ws.Range("J7").Value = 29.25
$ws.Range("K7").Value = 49.166668
$ws.Range("L7").Value = 29.25
$ws.Range("M7").Value = 63.833332
$ws.Range("N7").Value = -255.25
$ws.Range("H86").Value = 7417853.5
$ws.Range("I86").Value = 9536097
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 9536097
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -9534974
$ws.Range("N86").Value = -6246
$ws.Range("H89").Value = 7417853.5
$ws.Range("I89").Value = 9536097
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 47680485
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -47674869
$ws.Range("N89").Value = -31232
$ws.Range("H99").Value = 7712.274
$ws.Range("I99").Value = 8558.75
$ws.Range("K99").Value = 8558.75
$ws.Range("M99").Value = -7060.75
$ws.Range("H126").Value = 7712.274
$ws.Range("I126").Value = 8558.75
$ws.Range("K126").Value = 25676.25
$ws.Range("M126").Value = -23206.25

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 325.77777
$ws.Range("I33").Value = 97
$ws.Range("K33").Value = 582
$ws.Range("M33").Value = -299
$ws.Range("H63").Value = 5998.25
$ws.Range("I63").Value = 1997
$ws.Range("K63").Value = 5991
$ws.Range("M63").Value = -5242
$ws.Range("H66").Value = 5998.25
$ws.Range("I66").Value = 1997
$ws.Range("K66").Value = 17973
$ws.Range("M66").Value = -14229
$ws.Range("H75").Value = 387269.7
$ws.Range("I75").Value = 1877
$ws.Range("J75").Value = 502887.5
$ws.Range("K75").Value = 5631
$ws.Range("L75").Value = 1508662.5
$ws.Range("M75").Value = -4633
$ws.Range("N75").Value = -1510658.5
$ws.Range("H78").Value = 387269.7
$ws.Range("I78").Value = 1877
$ws.Range("J78").Value = 502887.5
$ws.Range("K78").Value = 16893
$ws.Range("L78").Value = 4525987.5
$ws.Range("M78").Value = -11901
$ws.Range("N78").Value = -4535971.5
$ws.Range("H98").Value = 519
$ws.Range("I98").Value = 446.66666
$ws.Range("K98").Value = 1339.99998
$ws.Range("M98").Value = 158.0000199999999
$ws.Range("H114").Value = 923.1667
$ws.Range("I114").Value = 1046.1428
$ws.Range("J114").Value = 844.9091
$ws.Range("K114").Value = 3138.4284
$ws.Range("L114").Value = 2534.7273
$ws.Range("M114").Value = 115.5715999999998
$ws.Range("N114").Value = -9042.7273
$ws.Range("H137").Value = 9454.454
$ws.Range("I137").Value = 3863
$ws.Range("J137").Value = 15045.909
$ws.Range("K137").Value = 11589
$ws.Range("L137").Value = 45137.727
$ws.Range("M137").Value = -6489
$ws.Range("N137").Value = -55337.727

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3228.7273
$ws.Range("I122").Value = 2501.7
$ws.Range("K122").Value = 7505.099999999999
$ws.Range("M122").Value = -5055.099999999999
$ws.Range("H126").Value = 7453.3687
$ws.Range("J126").Value = 9997
$ws.Range("L126").Value = 29991
$ws.Range("N126").Value = -34931

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2785.1667
$ws.Range("I7").Value = 2800
$ws.Range("K7").Value = 2800
$ws.Range("M7").Value = -2688
$ws.Range("H16").Value = 3475009.2
$ws.Range("I16").Value = 3908326
$ws.Range("J16").Value = 8474.5
$ws.Range("K16").Value = 3908326
$ws.Range("L16").Value = 8474.5
$ws.Range("M16").Value = -3908156
$ws.Range("N16").Value = -8814.5
$ws.Range("H40").Value = 4647.857
$ws.Range("I40").Value = 4644.5713
$ws.Range("J40").Value = 4654.4287
$ws.Range("K40").Value = 4644.5713
$ws.Range("L40").Value = 4654.4287
$ws.Range("M40").Value = -4508.5713
$ws.Range("N40").Value = -4926.4287
$ws.Range("H46").Value = 1780.3055
$ws.Range("I46").Value = 1341.591
$ws.Range("J46").Value = 2469.7144
$ws.Range("K46").Value = 1341.591
$ws.Range("L46").Value = 2469.7144
$ws.Range("M46").Value = -1153.591
$ws.Range("N46").Value = -2845.7144
$ws.Range("H61").Value = 7214.8667
$ws.Range("I61").Value = 7214.8667
$ws.Range("K61").Value = 7214.8667
$ws.Range("M61").Value = -7012.8667
$ws.Range("H113").Value = 7214.8667
$ws.Range("I113").Value = 7214.8667
$ws.Range("K113").Value = 7214.8667
$ws.Range("M113").Value = -5044.8667
$ws.Range("H126").Value = 2785.1667
$ws.Range("I126").Value = 2800
$ws.Range("K126").Value = 8400
$ws.Range("M126").Value = -5930
$ws.Range("H136").Value = 12089.02
$ws.Range("I136").Value = 12563.473
$ws.Range("K136").Value = 37690.419
$ws.Range("M136").Value = -35140.419
$ws.Range("H139").Value = 120000
$ws.Range("J139").Value = 120000
$ws.Range("L139").Value = 120000
$ws.Range("N139").Value = -130280

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 53073.1
$ws.Range("I100").Value = 666
$ws.Range("J100").Value = 88011.164
$ws.Range("K100").Value = 1332
$ws.Range("L100").Value = 176022.328
$ws.Range("M100").Value = -791
$ws.Range("N100").Value = -177104.328
$ws.Range("H126").Value = 5389.6113
$ws.Range("I126").Value = 5040.8
$ws.Range("K126").Value = 15122.4
$ws.Range("M126").Value = -12652.4
